# Atualização de bases das ligas, do dia: 15-06-2024 às 21:10
# Swap the data (columns B..AD) between paired rows, leaving column A (row index) untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstCol = 2   # column B
$lastCol  = 30  # column AD

$rowPairs = @(
    @(20, 21),
    @(28, 29),
    @(54, 55),
    @(56, 57)
)

foreach ($pair in $rowPairs) {
    $r1 = $pair[0]
    $r2 = $pair[1]

    # Read all values from both rows first
    $row1Vals = @{}
    $row2Vals = @{}
    for ($c = $firstCol; $c -le $lastCol; $c++) {
        $row1Vals[$c] = $ws.Cells.Item($r1, $c).Value2
        $row2Vals[$c] = $ws.Cells.Item($r2, $c).Value2
    }

    # Write swapped values back
    for ($c = $firstCol; $c -le $lastCol; $c++) {
        $ws.Cells.Item($r1, $c).Value2 = $row2Vals[$c]
        $ws.Cells.Item($r2, $c).Value2 = $row1Vals[$c]
    }
}
